# Update the "想去人数" (want-to-go count) column (F) for a handful of
# events on both the "展览" and "全部类型" sheets to reflect newly
# generated data (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F-column value, keyed by the values that need to change.
$updates = @{
    4  = 42
    6  = 148
    8  = 61
    10 = 5248
    11 = 4804
    16 = 185
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
